# This script updates the "Förändrad" (column C) date stamp for every data
# row from 45184 to 45186, and adds a friendly display-text second argument
# (the report's "Beteckning", column A) to every HYPERLINK() formula found
# in columns S, T, U, V, W, X and Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data.
$ur = $ws.UsedRange
$firstRow = $ur.Row
if ($firstRow -lt 1) { $firstRow = 1 }
$lastRow = $ur.Row + $ur.Rows.Count - 1

# Data starts on row 2 (row 1 is the header row).
$startRow = 2

# Columns that may contain HYPERLINK(...) formulas that need the
# display-text argument added.
$hyperlinkCols = 19, 20, 21, 22, 23, 24, 25   # S, T, U, V, W, X, Y

for ($r = $startRow; $r -le $lastRow; $r++) {

    # Column A holds the "Beteckning" identifier used as the hyperlink's
    # friendly display text.
    $beteckning = $ws.Cells.Item($r, 1).Value2

    # Column C ("Förändrad") -> bump date serial from 45184 to 45186.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like '*HYPERLINK(*' -and $f -notlike '*,*') {
                $trimmed = $f.TrimEnd()
                if ($trimmed.EndsWith(')')) {
                    $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $beteckning + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
